# Updated Master data as per 16th May Refresh
# Appends 3 new rows (34-36) of test data to the reg_center_user_h sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(10005, 110033),
    @(10005, 110034),
    @(10005, 110035)
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $regcntrId = $newRows[$i][0]
    $usrId = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $regcntrId
    $ws.Cells.Item($r, 2).Value = $usrId
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Match the post-edit selection state captured in the workbook.
$ws.Range("A37:XFD1048576").Select()
